$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (subject numbers) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) data values
$ws.Range("B2").Value = 15.319071850583157
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 20.262778902774098
$ws.Range("E2").ClearContents()

# Row 3 (STR) data values
$ws.Range("B3").Value = 13.361801302203441
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 15.925921680975783
$ws.Range("E3").Value = -13.535893596395896

# Update selection to reflect the edited range, matching the saved file
$ws.Range("B1:E3").Select()
